$d = $word.ActiveDocument

$found = $d.Content.Find.Execute(": All black suits +2", $true, $false, $false, $false, $false,
                         $true, 1, $false, ": All black suits +2 ", 2)

$range = $d.Content
$range.Find.Execute(": All black suits +2 ", $true, $false, $false, $false, $false,
                     $true, 1, $false, "", 0)

$insertRange = $d.Content
$insertRange.Find.Execute(": All black suits +2 ")
$insertRange.Collapse(0)
$newRange = $insertRange.InsertAfter("IMPLEMENTED")
$newRange.Font.Bold = $true
$newRange.Font.Size = 16
$newRange.LanguageID = 1033
